# Add a "Slovakia" market sheet, cloned from the existing "Portugal" sheet,
# appended as the new last/active tab — mirrors Excel's own
# Worksheets("Portugal").Copy After:=<last sheet> workflow.

$wb = $excel.ActiveWorkbook

$portugal = $wb.Worksheets.Item("Portugal")

# Select the full sheet on the source before copying, matching the
# "select all" state the source sheet is left in after the tab is cloned.
$portugal.Range("A1:XFD1048576").Select() | Out-Null

# Copy "Portugal" to the end of the workbook; this creates the new sheet,
# makes it the active tab, and updates bookViews/activeTab accordingly.
$portugal.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Slovakia"

# Slovakia's repeater list doesn't include the P32AR/P32DR rows that
# Portugal has (rows 16 & 17) — remove them, shifting the trailing
# Wg/Repeaters rows up.
$newSheet.Rows.Item(16).Delete() | Out-Null
$newSheet.Rows.Item(16).Delete() | Out-Null

# Market-specific user story reference.
$newSheet.Range("B4").Value = "NGC-2930/T3178"

# Leave the new sheet's selection on B2, as a freshly-copied tab would show.
$newSheet.Range("B2").Select() | Out-Null
